$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.609.61"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "2.641.16"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'525.34"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.77%  "
$ws.Range("D6").Value = "'153.78"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.577"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.07%  "
$ws.Range("D9").Value = "'6.47"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.46%  "
$ws.Range("E10").Value = "  +3.15%  "
$ws.Range("D11").Value = "'0.348"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("D13").Value = "3.102.28"
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("D14").Value = "60.577.15"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").Value = "'21.82"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "2.645.19"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("D19").Value = "'351.12"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.73%  "
$ws.Range("D20").Value = "'10.59"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("D22").Value = "'0.996"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").Value = "'61.15"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("D27").Value = "0.0₃0845"
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").Value = "'7.23"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.17%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "'6.15"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.90%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.61"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.76%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'19.28"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("D33").Value = "'149.66"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.98%  "
$ws.Range("D34").Value = "'4.07"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.87%  "
$ws.Range("E35").Value = "  -1.01%  "
$ws.Range("D36").Value = "'0.899"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.35%  "
$ws.Range("D37").Value = "'0.886"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.73%  "
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("D39").Value = "'305.12"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.33%  "
$ws.Range("D40").Value = "'1.47"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.76%  "
$ws.Range("E41").Value = "  -0.58%  "
$ws.Range("D42").Value = "'0.635"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.81%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").Value = "'0.998"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "'0.0557"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'19.93"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("D47").Value = "'0.0238"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("E48").Value = "  -3.13%  "
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("D50").Value = "'19.01"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("D51").Value = "1.978.99"
$ws.Range("E51").Value = "  -0.90%  "

Write-Output "done"
